# Apply the "Add canvas, transition. machine learning seconde commit" edit
# to TestData.xlsx / Sheet1.
#
# Summary of the change (per the OOXML diff):
#   1. A bunch of previously-zero cells in columns B:F (rows 2-21) become 1
#      (B3 becomes 2 specifically).
#   2. The sheet view scrolls down and the selection becomes the whole
#      used range A1:G21 (active cell around A9).
#   3. A <pageSetup> (paperSize=9/A4, portrait) is added to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Cell value updates (all former 0s -> 1, except B3 -> 2) ---------

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

# Row 3
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Rows 4-8 (column C only)
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1

# Rows 10-12
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F12").Value = 1

# Rows 13-15
$ws.Range("B13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = 1

$ws.Range("B14").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("F14").Value = 1

$ws.Range("B15").Value = 1
$ws.Range("D15").Value = 1

# Rows 18-21
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("D21").Value = 1

# --- 2. View / selection: scroll down and select the whole table --------

$excel.ActiveWindow.ScrollRow = 9
$ws.Range("A1:G21").Select()

# --- 3. Page setup: paper size 9 (A4), portrait orientation --------------

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = "xlPortrait"
